# Update scripts with new TPM values: the LR-pair matrix (ECs/FAPs/MuSCs as
# sending clusters, FAPs/MuSCs as target clusters, Slurp1 -> Chrna7) was
# recomputed against a new TPM dataset. Rewrite all data rows (2-7) with the
# refreshed cluster labels and expression/specificity figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slurp1"
$ws.Range("C2").Value = "Chrna7"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.033624
$ws.Range("H2").Value = 0.100872
$ws.Range("I2").Value = 0.03079659537519841
$ws.Range("J2").Value = 0.03079659537519841
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2159606666666667
$ws.Range("N2").Value = 0.647882
$ws.Range("O2").Value = 0.09879736522873545
$ws.Range("P2").Value = 0.09879736522873543
$ws.Range("Q2").Value = 0.007261461456
$ws.Range("R2").Value = 0.065353153104
$ws.Range("S2").Value = 0.003042622481085062
$ws.Range("T2").Value = 0.003042622481085062

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slurp1"
$ws.Range("C3").Value = "Chrna7"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.033624
$ws.Range("H3").Value = 0.100872
$ws.Range("I3").Value = 0.03079659537519841
$ws.Range("J3").Value = 0.03079659537519841
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.969934333333333
$ws.Range("N3").Value = 5.909803
$ws.Range("O3").Value = 0.9012026347712646
$ws.Range("P3").Value = 0.9012026347712645
$ws.Range("Q3").Value = 0.066237072024
$ws.Range("R3").Value = 0.596133648216
$ws.Range("S3").Value = 0.02775397289411335
$ws.Range("T3").Value = 0.02775397289411335

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Slurp1"
$ws.Range("C4").Value = "Chrna7"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9671793333333333
$ws.Range("H4").Value = 2.901538
$ws.Range("I4").Value = 0.8858503028765409
$ws.Range("J4").Value = 0.8858503028765411
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2159606666666667
$ws.Range("N4").Value = 0.647882
$ws.Range("O4").Value = 0.09879736522873545
$ws.Range("P4").Value = 0.09879736522873543
$ws.Range("Q4").Value = 0.2088726936128889
$ws.Range("R4").Value = 1.879854242516
$ws.Range("S4").Value = 0.08751967591127952
$ws.Range("T4").Value = 0.08751967591127953

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slurp1"
$ws.Range("C5").Value = "Chrna7"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9671793333333333
$ws.Range("H5").Value = 2.901538
$ws.Range("I5").Value = 0.8858503028765409
$ws.Range("J5").Value = 0.8858503028765411
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.969934333333333
$ws.Range("N5").Value = 5.909803
$ws.Range("O5").Value = 0.9012026347712646
$ws.Range("P5").Value = 0.9012026347712645
$ws.Range("Q5").Value = 1.905279775223778
$ws.Range("R5").Value = 17.147517977014
$ws.Range("S5").Value = 0.7983306269652614
$ws.Range("T5").Value = 0.7983306269652615

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Slurp1"
$ws.Range("C6").Value = "Chrna7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.09100566666666665
$ws.Range("H6").Value = 0.273017
$ws.Range("I6").Value = 0.08335310174826059
$ws.Range("J6").Value = 0.0833531017482606
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2159606666666667
$ws.Range("N6").Value = 0.647882
$ws.Range("O6").Value = 0.09879736522873545
$ws.Range("P6").Value = 0.09879736522873543
$ws.Range("Q6").Value = 0.01965364444377777
$ws.Range("R6").Value = 0.1768827999939999
$ws.Range("S6").Value = 0.008235066836370849
$ws.Range("T6").Value = 0.008235066836370849

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Slurp1"
$ws.Range("C7").Value = "Chrna7"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.09100566666666665
$ws.Range("H7").Value = 0.273017
$ws.Range("I7").Value = 0.08335310174826059
$ws.Range("J7").Value = 0.0833531017482606
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.969934333333333
$ws.Range("N7").Value = 5.909803
$ws.Range("O7").Value = 0.9012026347712646
$ws.Range("P7").Value = 0.9012026347712645
$ws.Range("Q7").Value = 0.1792751872945555
$ws.Range("R7").Value = 1.613476685651
$ws.Range("S7").Value = 0.07511803491188974
$ws.Range("T7").Value = 0.07511803491188974
